$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.2054769113617
$ws.Range("C2").Value = 0.3299325156323505
$ws.Range("D2").Value = 0.01388337651986404
$ws.Range("E2").Value = 0.4245925918300344
$ws.Range("F2").Value = 0.441698669717745
$ws.Range("O2").Value = 1.383927038744844

$ws.Range("B3").Value = 1.056831324844325
$ws.Range("C3").Value = 0.2904986091814408
$ws.Range("D3").Value = 0.01236594077360564
$ws.Range("E3").Value = 0.37026579236597
$ws.Range("F3").Value = 0.4375669543272593
$ws.Range("O3").Value = 1.384955041402407

$ws.Range("B4").Value = 0.9653392133317311
$ws.Range("C4").Value = 0.2661785677906607
$ws.Range("D4").Value = 0.01143135250730154
$ws.Range("E4").Value = 0.3370109357875037
$ws.Range("F4").Value = 0.4355725835690976
$ws.Range("O4").Value = 1.387419968358103

$ws.Range("B5").Value = 0.9280008519348257
$ws.Range("C5").Value = 0.2562414314373882
$ws.Range("D5").Value = 0.01104980545873246
$ws.Range("E5").Value = 0.3234828543049986
$ws.Range("F5").Value = 0.4348954029591994
$ws.Range("O5").Value = 1.388882426217862

$ws.Range("B6").Value = 0.9217975943939791
$ws.Range("C6").Value = 0.2545897906098844
$ws.Range("D6").Value = 0.01098640875816415
$ws.Range("E6").Value = 0.3212378872433845
$ws.Range("F6").Value = 0.4347911182924733
$ws.Range("O6").Value = 1.38915283563486

$ws.Range("B7").Value = 0.9648358740913636
$ws.Range("C7").Value = 0.2660446588068339
$ws.Range("D7").Value = 0.01142620960523999
$ws.Range("E7").Value = 0.3368283987874321
$ws.Range("F7").Value = 0.435562903239834
$ws.Range("O7").Value = 1.387437841387509

$ws.Range("B8").Value = 1.154271079244324
$ws.Range("C8").Value = 0.3163582940311471
$ws.Range("D8").Value = 0.01336077815361847
$ws.Range("E8").Value = 0.4058382475447218
$ws.Range("F8").Value = 0.440160922482022
$ws.Range("O8").Value = 1.383899017092546

$ws.Range("B9").Value = 1.523935298033621
$ws.Range("C9").Value = 0.4141536311091727
$ws.Range("D9").Value = 0.01713059949575069
$ws.Range("E9").Value = 0.5420775926543655
$ws.Range("F9").Value = 0.453523150812174
$ws.Range("O9").Value = 1.391652447194474

$ws.Range("B10").Value = 1.794384906139726
$ws.Range("C10").Value = 0.4854586294006822
$ws.Range("D10").Value = 0.01988466917678267
$ws.Range("E10").Value = 0.6428785719488133
$ws.Range("F10").Value = 0.4660478875119622
$ws.Range("O10").Value = 1.406507688373466

$ws.Range("B11").Value = 1.917166366813888
$ws.Range("C11").Value = 0.5177761598071697
$ws.Range("D11").Value = 0.02113398204294725
$ws.Range("E11").Value = 0.6889221478109562
$ws.Range("F11").Value = 0.4723463114223563
$ws.Range("O11").Value = 1.415298141285945

$ws.Range("B12").Value = 1.963623957007314
$ws.Range("C12").Value = 0.5299964307634468
$ws.Range("D12").Value = 0.02160653608757457
$ws.Range("E12").Value = 0.7063873339732254
$ws.Range("F12").Value = 0.4748187584742567
$ws.Range("O12").Value = 1.418922657074035

$ws.Range("B13").Value = 1.953620156600778
$ws.Range("C13").Value = 0.5273653707092194
$ws.Range("D13").Value = 0.0215047872627494
$ws.Range("E13").Value = 0.7026245446588746
$ws.Range("F13").Value = 0.4742823714906024
$ws.Range("O13").Value = 1.418128843461943

$ws.Range("B14").Value = 1.920989212423251
$ws.Range("C14").Value = 0.5187818859170648
$ws.Range("D14").Value = 0.02117287019112268
$ws.Range("E14").Value = 0.6903584147401745
$ws.Range("F14").Value = 0.4725479645388617
$ws.Range("O14").Value = 1.415590385022313

$ws.Range("B15").Value = 1.900996943457415
$ws.Range("C15").Value = 0.5135219371092035
$ws.Range("D15").Value = 0.02096949113956725
$ws.Range("E15").Value = 0.6828489657860644
$ws.Range("F15").Value = 0.4714969972042695
$ws.Range("O15").Value = 1.414074125701575

$ws.Range("B16").Value = 1.786355740578244
$ws.Range("C16").Value = 0.4833441559972584
$ws.Range("D16").Value = 0.01980295038744373
$ws.Range("E16").Value = 0.6398735132023745
$ws.Range("F16").Value = 0.4656484441841258
$ws.Range("O16").Value = 1.40597441066285

$ws.Range("B17").Value = 1.7159626214073
$ws.Range("C17").Value = 0.4648001082296673
$ws.Range("D17").Value = 0.0190863929290046
$ws.Range("E17").Value = 0.6135594577242642
$ws.Range("F17").Value = 0.4622151284550426
$ws.Range("O17").Value = 1.401528585528638

$ws.Range("B18").Value = 1.67545115994136
$ws.Range("C18").Value = 0.4541228461810647
$ws.Range("D18").Value = 0.01867391755196479
$ws.Range("E18").Value = 0.5984419285571221
$ws.Range("F18").Value = 0.4602968736919166
$ws.Range("O18").Value = 1.399162587266261

$ws.Range("B19").Value = 1.661730736876791
$ws.Range("C19").Value = 0.4505057976432454
$ws.Range("D19").Value = 0.01853420453561938
$ws.Range("E19").Value = 0.593326348084986
$ws.Range("F19").Value = 0.4596570596562657
$ws.Range("O19").Value = 1.398394221217472

$ws.Range("B20").Value = 1.72345850135747
$ws.Range("C20").Value = 0.4667753182363299
$ws.Range("D20").Value = 0.01916270609149962
$ws.Range("E20").Value = 0.6163587933429682
$ws.Range("F20").Value = 0.4625747564943907
$ws.Range("O20").Value = 1.401982045758189

$ws.Range("B21").Value = 1.930574729968157
$ws.Range("C21").Value = 0.521303547098853
$ws.Range("D21").Value = 0.02127037696963896
$ws.Range("E21").Value = 0.6939604554036123
$ws.Range("F21").Value = 0.4730550229627255
$ws.Range("O21").Value = 1.416327937285843

$ws.Range("B22").Value = 2.065720332583794
$ws.Range("C22").Value = 0.5568376430126136
$ws.Range("D22").Value = 0.02264473974224757
$ws.Range("E22").Value = 0.7448505902177658
$ws.Range("F22").Value = 0.4804141968492672
$ws.Range("O22").Value = 1.42742934729705

$ws.Range("B23").Value = 1.993610907615164
$ws.Range("C23").Value = 0.5378820360190844
$ws.Range("D23").Value = 0.02191151061945362
$ws.Range("E23").Value = 0.7176729576034973
$ws.Range("F23").Value = 0.4764395078873349
$ws.Range("O23").Value = 1.42134528323092

$ws.Range("B24").Value = 1.7200697432487
$ws.Range("C24").Value = 0.4658823757277446
$ws.Range("D24").Value = 0.01912820651954661
$ws.Range("E24").Value = 0.6150931803568795
$ws.Range("F24").Value = 0.4624119955302888
$ws.Range("O24").Value = 1.401776444925162

$ws.Range("B25").Value = 1.424129652389411
$ws.Range("C25").Value = 0.3877921022066744
$ws.Range("D25").Value = 0.01611343166023005
$ws.Range("E25").Value = 0.4245925918300344
$ws.Range("F25").Value = 0.4494370706573534
$ws.Range("O25").Value = 1.387961149817301

